$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">1. </w:t>
      </w:r>
      <w:r>
        <w:t>In your opinion, out of the given test cases, list the easiest test case(s). Why? [1 pt]</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test1.micro – It only tests that a symbol table is generated without having any code blocks</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test5.micro – Only one symbol table</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test14.micro – it tests comments, so there isn’t a lot being added to the symbol table</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test16.micro – it has two empty symbol tables</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">2. </w:t>
      </w:r>
      <w:r>
        <w:t>In your opinion, out of the given test cases, list the hardest test cases(s). Why? [1 pt]</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>test9.micro – error handling and symbol table dropping</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test20.micro – global string declaration with 7 symbol tables, one with a custom name</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>test21.micro – error handling in a separate method block</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">3. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">List all </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterX</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitX</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> functions that you think need to be implemented. (</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>hint</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">: check the auto-generated </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>LittleBaseListener</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> class) [4 pts]</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enter_Program</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exit_Program</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterFunc</w:t>
      </w:r>
      <w:r>
        <w:t>_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>e</w:t>
      </w:r>
      <w:r>
        <w:t>xitFunc</w:t>
      </w:r>
      <w:r>
        <w:t>_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterP</w:t>
      </w:r>
      <w:r>
        <w:t>gm_body</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitP</w:t>
      </w:r>
      <w:r>
        <w:t>gm_body</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterIf_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitIf_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterElse_part</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitElse_part</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterFor_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitFor_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterWhile_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitWhile_stmt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>enterString_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitString_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>e</w:t>
      </w:r>
      <w:r>
        <w:t>nterVar_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>exitVar_decl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">4. </w:t>
      </w:r>
      <w:r>
        <w:t>In your opinion, which Java data structure(s) should be used to implement the symbol table(s)? Why? [2 pts]</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Hash tables provide</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> a consistently high-performance solution due to their ability to insert and retrieve symbols in constant time on average. This means symbol tables can grow or shrink with each scope and performance will not be impacted.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">5. </w:t>
      </w:r>
      <w:r>
        <w:t>You also have the option of using a Visitor (instead of the Listener). What is the main difference between a Listener and a Visitor? (</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>hint</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>: read the above excerpt) [2 pts]</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>The main difference between the two is how the methods are called</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">n a visitor, methods have to give their children specific </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>calls</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> or they won’t be visited;</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> however in a Listener</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, the methods are called by an </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Antlr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>-provided walker object</w:t>
      </w:r>
      <w:r>
        <w:t>. Visitors are useful in situations where you may need to return data from nodes, or if you need to visit several nodes before performing an operation to change the parse tree.</w:t>
      </w:r>
    </w:p>
    <w:sectPr w:rsidR="00C17632">
      <w:pgSz w:w="12240" w:h="15840"/>
      <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/>
      <w:cols w:space="720"/>
      <w:docGrid w:linePitch="360"/>
    </w:sectPr>
  </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
